$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store numeric-looking / percent-looking
# values as plain text. Force the target range to Text format before writing
# so Excel does not silently reinterpret strings like "304.43" or "0.03%" as
# numbers, then restore the default "Normal" style so the sheet keeps its
# original (unstyled) look for these data cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "304.43"
$ws.Range("E2").Value = "0.03%"
$ws.Range("D3").Value = "36.16"
$ws.Range("E3").Value = "-2.16%"
$ws.Range("D4").Value = "5.035"
$ws.Range("E4").Value = "-0.08%"
$ws.Range("D5").Value = "0.07858"
$ws.Range("E5").Value = "-0.48%"
$ws.Range("D6").Value = "2.132"
$ws.Range("E6").Value = "-2.89%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "7.956"
$ws.Range("E7").Value = "-0.67%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9216"
$ws.Range("E8").Value = "-0.73%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.09532"
$ws.Range("E9").Value = "-3.91%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1851"
$ws.Range("E10").Value = "-1.54%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08591"
$ws.Range("E11").Value = "-0.94%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03611"
$ws.Range("E12").Value = "0.17%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09920"
$ws.Range("E13").Value = "-0.37%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001435"
$ws.Range("E14").Value = "-3.67%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005670"
$ws.Range("E15").Value = "0.38%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.470"
$ws.Range("E16").Value = "0.36%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.141"
$ws.Range("E17").Value = "2.63%"
$ws.Range("D18").Value = "2.752"
$ws.Range("E18").Value = "10.85%"
$ws.Range("D20").Value = "0.1337"
$ws.Range("E20").Value = "0.88%"
$ws.Range("D21").Value = "5.169"
$ws.Range("E21").Value = "8.34%"
$ws.Range("E22").Value = "2.39%"
$ws.Range("D23").Value = "0.04578"
$ws.Range("E23").Value = "-0.41%"
$ws.Range("E24").Value = "-1.30%"
$ws.Range("D25").Value = "0.004794"
$ws.Range("E25").Value = "-8.53%"
$ws.Range("E26").Value = "-7.03%"
$ws.Range("D27").Value = "0.0004751"
$ws.Range("E27").Value = "75.01%"
$ws.Range("D39").Value = "0.01862"
$ws.Range("E39").Value = "1.67%"
$ws.Range("D40").Value = "0.04707"
$ws.Range("E40").Value = "-1.33%"
$ws.Range("E41").Value = "-2.29%"
$ws.Range("D42").Value = "0.1386"
$ws.Range("E42").Value = "-1.97%"
$ws.Range("D43").Value = "0.007728"
$ws.Range("E43").Value = "2.58%"
$ws.Range("D44").Value = "0.002231"
$ws.Range("E44").Value = "1.95%"
$ws.Range("D45").Value = "0.01119"
$ws.Range("E45").Value = "7.31%"
$ws.Range("D46").Value = "0.00006358"
$ws.Range("E46").Value = "0.95%"
$ws.Range("E47").Value = "0.16%"
$ws.Range("E48").Value = "0.29%"
$ws.Range("D49").Value = "51.58"
$ws.Range("E49").Value = "41.96%"
$ws.Range("D50").Value = "0.001901"
$ws.Range("E50").Value = "-29.25%"
$ws.Range("D51").Value = "0.00002101"
$ws.Range("E51").Value = "0.16%"

$ws.Range("D2:E51").Style = "Normal"
